$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 139
$ws.Range("I6").Value = 139
$ws.Range("K6").Value = 417
$ws.Range("M6").Value = -305

$ws.Range("H33").Value = 702.55554
$ws.Range("I33").Value = 786.2
$ws.Range("K33").Value = 786.2
$ws.Range("M33").Value = -557.2

$ws.Range("H42").Value = 158.72728
$ws.Range("I42").Value = 32.25
$ws.Range("J42").Value = 310.5
$ws.Range("K42").Value = 96.75
$ws.Range("L42").Value = 931.5
$ws.Range("M42").Value = 133.25
$ws.Range("N42").Value = -1391.5

$ws.Range("H74").Value = 4786.923
$ws.Range("I74").Value = 1949.4
$ws.Range("J74").Value = 6560.375
$ws.Range("K74").Value = 1949.4
$ws.Range("L74").Value = 6560.375
$ws.Range("M74").Value = -1013.4
$ws.Range("N74").Value = -8432.375

$ws.Range("H77").Value = 4786.923
$ws.Range("I77").Value = 1949.4
$ws.Range("J77").Value = 6560.375
$ws.Range("K77").Value = 9747
$ws.Range("L77").Value = 32801.875
$ws.Range("M77").Value = -5067
$ws.Range("N77").Value = -42161.875

$ws.Range("H116").Value = 6240.773
$ws.Range("I116").Value = 4726.636
$ws.Range("J116").Value = 7754.909
$ws.Range("K116").Value = 4726.636
$ws.Range("L116").Value = 7754.909
$ws.Range("M116").Value = -1284.636
$ws.Range("N116").Value = -14638.909

$ws.Range("H125").Value = 8549779
$ws.Range("I125").Value = 1336.8
$ws.Range("K125").Value = 12031.2
$ws.Range("M125").Value = -9571.199999999999

$ws.Range("H137").Value = 51509.195
$ws.Range("I137").Value = 59053.547
$ws.Range("K137").Value = 177160.641
$ws.Range("M137").Value = -174610.641

$ws.Range("H141").Value = 10840.615
$ws.Range("I141").Value = 10840.615
$ws.Range("K141").Value = 32521.845
$ws.Range("M141").Value = -27341.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12601.704
$ws.Range("I32").Value = 7799.794
$ws.Range("J32").Value = 20764.95
$ws.Range("K32").Value = 7799.794
$ws.Range("L32").Value = 20764.95
$ws.Range("M32").Value = -7512.794
$ws.Range("N32").Value = -21338.95

$ws.Range("H61").Value = 3623.4443
$ws.Range("I61").Value = 3259.4285
$ws.Range("K61").Value = 3259.4285
$ws.Range("M61").Value = -3047.4285

$ws.Range("H74").Value = 20164.652
$ws.Range("I74").Value = 1255.5807
$ws.Range("K74").Value = 1255.5807
$ws.Range("M74").Value = -381.5807

$ws.Range("H77").Value = 20164.652
$ws.Range("I77").Value = 1255.5807
$ws.Range("K77").Value = 6277.9035
$ws.Range("M77").Value = -1909.9035

$ws.Range("H136").Value = 3623.4443
$ws.Range("I136").Value = 3259.4285
$ws.Range("K136").Value = 9778.2855
$ws.Range("M136").Value = -7228.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 5000
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5450

$ws.Range("H67").Value = 5000
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6560

$ws.Range("H105").Value = 8929528
$ws.Range("I105").Value = 8929528
$ws.Range("K105").Value = 8929528
$ws.Range("M105").Value = -8927781

$ws.Range("H107").Value = 2749853
$ws.Range("I107").Value = 4203571
$ws.Range("J107").Value = 3940.889
$ws.Range("K107").Value = 4203571
$ws.Range("L107").Value = 3940.889
$ws.Range("M107").Value = -4201651
$ws.Range("N107").Value = -7780.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1517.6666
$ws.Range("I16").Value = 1290.4615
$ws.Range("J16").Value = 2994.5
$ws.Range("K16").Value = 1290.4615
$ws.Range("L16").Value = 2994.5
$ws.Range("M16").Value = -1003.4615
$ws.Range("N16").Value = -3568.5

$ws.Range("H31").Value = 20454.666
$ws.Range("I31").Value = 2520.4
$ws.Range("K31").Value = 2520.4
$ws.Range("M31").Value = -2225.4

$ws.Range("H34").Value = 20454.666
$ws.Range("I34").Value = 2520.4
$ws.Range("K34").Value = 2520.4
$ws.Range("M34").Value = -2318.4

$ws.Range("H58").Value = 4463.5713
$ws.Range("I58").Value = 5765.4287
$ws.Range("J58").Value = 3161.7144
$ws.Range("K58").Value = 5765.4287
$ws.Range("L58").Value = 3161.7144
$ws.Range("M58").Value = -5562.4287
$ws.Range("N58").Value = -3567.7144

$ws.Range("H113").Value = 1517.6666
$ws.Range("I113").Value = 1290.4615
$ws.Range("J113").Value = 2994.5
$ws.Range("K113").Value = 1290.4615
$ws.Range("L113").Value = 2994.5
$ws.Range("M113").Value = 879.5385000000001
$ws.Range("N113").Value = -7334.5

$ws.Range("H122").Value = 2168.3462
$ws.Range("I122").Value = 2065.7083
$ws.Range("J122").Value = 3400
$ws.Range("K122").Value = 6197.124899999999
$ws.Range("L122").Value = 10200
$ws.Range("M122").Value = -3747.124899999999
$ws.Range("N122").Value = -15100

$ws.Range("H136").Value = 4463.5713
$ws.Range("I136").Value = 5765.4287
$ws.Range("J136").Value = 3161.7144
$ws.Range("K136").Value = 17296.2861
$ws.Range("L136").Value = 9485.143199999999
$ws.Range("M136").Value = -14746.2861
$ws.Range("N136").Value = -14585.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 250469.25
$ws.Range("I9").Value = 333700
$ws.Range("J9").Value = 777
$ws.Range("K9").Value = 1001100
$ws.Range("L9").Value = 2331
$ws.Range("M9").Value = -1000876
$ws.Range("N9").Value = -2779

$ws.Range("H36").Value = 370.33334
$ws.Range("I36").Value = 370.33334
$ws.Range("K36").Value = 1111.00002
$ws.Range("M36").Value = -942.0000199999999

$ws.Range("H80").Value = 1999.5
$ws.Range("J80").Value = 1999.5
$ws.Range("L80").Value = 5998.5
$ws.Range("N80").Value = -7870.5

$ws.Range("H83").Value = 1999.5
$ws.Range("J83").Value = 1999.5
$ws.Range("L83").Value = 17995.5
$ws.Range("N83").Value = -27355.5

$ws.Range("H92").Value = 664
$ws.Range("I92").Value = 218.11111
$ws.Range("K92").Value = 654.3333299999999
$ws.Range("M92").Value = 593.6666700000001

$ws.Range("H107").Value = 1886.3334
$ws.Range("J107").Value = 1038
$ws.Range("L107").Value = 3114
$ws.Range("N107").Value = -6954

$ws.Range("H121").Value = 1205.1538
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1205.1538
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3615.4614
$ws.Range("N121").Value = -6235.4614
$ws.Range("M121").ClearContents()

$ws.Range("H129").Value = 1496.3529
$ws.Range("I129").Value = 1105.4
$ws.Range("J129").Value = 1659.25
$ws.Range("K129").Value = 3316.2
$ws.Range("L129").Value = 4977.75
$ws.Range("M129").Value = 1683.8
$ws.Range("N129").Value = -14977.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5379289.5
$ws.Range("I113").Value = 11112669
$ws.Range("K113").Value = 11112669
$ws.Range("M113").Value = -11110499

$ws.Range("H122").Value = 308918.22
$ws.Range("I122").Value = 469957.94
$ws.Range("K122").Value = 1409873.82
$ws.Range("M122").Value = -1407423.82

$ws.Range("H132").Value = 3308
$ws.Range("I132").Value = 3082.9
$ws.Range("K132").Value = 9248.700000000001
$ws.Range("M132").Value = -6718.700000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 211114540
$ws.Range("I82").Value = 211114540
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 211114540
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -211114179
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 211114540
$ws.Range("I85").Value = 211114540
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 211114540
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -211113292
$ws.Range("N85").ClearContents()

$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27246

$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -86232

$ws.Range("H93").Value = 83359840
$ws.Range("I93").Value = 111113110
$ws.Range("J93").Value = 99999
$ws.Range("K93").Value = 111113110
$ws.Range("L93").Value = 99999
$ws.Range("M93").Value = -111111862
$ws.Range("N93").Value = -102495

$ws.Range("H132").Value = 8072.3213
$ws.Range("I132").Value = 8629.348
$ws.Range("K132").Value = 25888.044
$ws.Range("M132").Value = -23358.044

$ws.Range("H134").Value = 81357
$ws.Range("J134").Value = 81357
$ws.Range("L134").Value = 81357
$ws.Range("N134").Value = -91497

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 3344.1428
$ws.Range("J23").Value = 7299.6665
$ws.Range("L23").Value = 7299.6665
$ws.Range("N23").Value = -7757.6665

$ws.Range("H62").Value = 8717.9
$ws.Range("J62").Value = 9108.777
$ws.Range("L62").Value = 9108.777
$ws.Range("N62").Value = -10356.777

$ws.Range("H65").Value = 8717.9
$ws.Range("J65").Value = 9108.777
$ws.Range("L65").Value = 45543.885
$ws.Range("N65").Value = -51783.885

$ws.Range("H96").Value = 3048.88
$ws.Range("I96").Value = 2748.7368
$ws.Range("J96").Value = 3999.3333
$ws.Range("K96").Value = 2748.7368
$ws.Range("L96").Value = 3999.3333
$ws.Range("M96").Value = -1375.7368
$ws.Range("N96").Value = -6745.3333

$ws.Range("H109").Value = 66996.336
$ws.Range("J109").Value = 66996.336
$ws.Range("L109").Value = 66996.336
$ws.Range("N109").Value = -69770.336

$ws.Range("H136").Value = 1737.7142
$ws.Range("I136").Value = 1266.8334
$ws.Range("K136").Value = 3800.5002
$ws.Range("M136").Value = -1250.5002
